# Generate Report for Handoff
# Replaces the tracked source file (old GUID) with a newly generated one
# (new GUID), refreshes the handoff timestamps/xlf artifact names, and
# clears the stale handback bookkeeping (target/handback file + datetime)
# for the file that has just gone out for handoff again.

$wb = $excel.ActiveWorkbook

$oldId = "2568fbfd-1eb4-41a6-8c96-17490e131c4f"
$newId = "46f920be-dd7f-41e4-a422-67af402ee598"

$newMd      = "$newId.md"
$newMdPath  = "e2e\$newId.md"
$newZhCnXlf = "$newId.fbfe421b1826b01d6037d9048556b0093386ccf2.zh-cn.xlf"
$newDeDeXlf = "$newId.fbfe421b1826b01d6037d9048556b0093386ccf2.de-de.xlf"

$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = "2016-08-29 17:07:04"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = $newMdPath
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = "2016-08-29 17:06:56"
$wsZhCn.Range("I2").Value = "'"
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = "'"
$wsZhCn.Range("J2").Style = "Normal"
$wsZhCn.Range("K2").Value = $zeroDate

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsZhCn.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZhCn.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = "2016-08-29 17:07:04"
$wsDeDe.Range("I2").Value = "'"
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = "'"
$wsDeDe.Range("J2").Style = "Normal"
$wsDeDe.Range("K2").Value = $zeroDate

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDeDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDeDe.Columns.Item(10).ColumnWidth = 20.833333333333332
